$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Cells.Item(8, 5).Value = "فصل سوم منتهی به 1399/09"
$ws.Cells.Item(8, 6).Value = "فصل چهارم منتهی به 1399/12"
$ws.Cells.Item(8, 7).Value = "فصل اول منتهی به 1400/03"
$ws.Cells.Item(8, 8).Value = "فصل دوم منتهی به 1400/06"
$ws.Cells.Item(8, 9).Value = "فصل سوم منتهی به 1400/09"
$ws.Cells.Item(8, 10).Value = "فصل چهارم منتهی به 1400/12"
$ws.Cells.Item(8, 11).Value = "فصل اول منتهی به 1401/03"
$ws.Cells.Item(8, 12).Value = "فصل دوم منتهی به 1401/06"
$ws.Cells.Item(8, 13).Value = "فصل سوم منتهی به 1401/09"
$ws.Cells.Item(8, 14).Value = "فصل چهارم منتهی به 1401/12"

# Row 11
$ws.Cells.Item(11, 5).Value = 711833
$ws.Cells.Item(11, 6).Value = 1835712
$ws.Cells.Item(11, 7).Value = 441238
$ws.Cells.Item(11, 8).Value = 359613
$ws.Cells.Item(11, 9).Value = 307350
$ws.Cells.Item(11, 10).Value = 213137
$ws.Cells.Item(11, 11).Value = 365400
$ws.Cells.Item(11, 12).Value = 339706
$ws.Cells.Item(11, 13).Value = 397486
$ws.Cells.Item(11, 14).Value = 657401

# Row 12
$ws.Cells.Item(12, 5).Value = 156591
$ws.Cells.Item(12, 6).Value = 513572
$ws.Cells.Item(12, 7).Value = 56407
$ws.Cells.Item(12, 8).Value = 91604
$ws.Cells.Item(12, 9).Value = 141204
$ws.Cells.Item(12, 10).Value = 117159
$ws.Cells.Item(12, 11).Value = 76103
$ws.Cells.Item(12, 12).Value = 14866
$ws.Cells.Item(12, 13).Value = 72484
$ws.Cells.Item(12, 14).Value = 94922

# Row 13
$ws.Cells.Item(13, 8).Value = 0

# Row 33
$ws.Cells.Item(33, 10).Value = 22441475
$ws.Cells.Item(33, 11).Value = 1029600
$ws.Cells.Item(33, 12).Value = 3864079
$ws.Cells.Item(33, 13).Value = 1636911
$ws.Cells.Item(33, 14).Value = 1289874

# Row 34
$ws.Cells.Item(34, 5).Value = 266955
$ws.Cells.Item(34, 6).Value = 1000247
$ws.Cells.Item(34, 7).Value = 10881870
$ws.Cells.Item(34, 8).Value = 41655981
$ws.Cells.Item(34, 9).Value = 14430693
$ws.Cells.Item(34, 10).Value = "-"

# Row 35
$ws.Cells.Item(35, 5).Value = 1135379
$ws.Cells.Item(35, 6).Value = 3349531
$ws.Cells.Item(35, 7).Value = 11379515
$ws.Cells.Item(35, 8).Value = 42107198
$ws.Cells.Item(35, 9).Value = 14879247
$ws.Cells.Item(35, 10).Value = 22771771
$ws.Cells.Item(35, 11).Value = 1471103
$ws.Cells.Item(35, 12).Value = 4218651
$ws.Cells.Item(35, 13).Value = 2106881
$ws.Cells.Item(35, 14).Value = 2042197

# Row 39
$ws.Cells.Item(39, 5).Value = "فصل سوم منتهی به 1399/09"
$ws.Cells.Item(39, 6).Value = "فصل چهارم منتهی به 1399/12"
$ws.Cells.Item(39, 7).Value = "فصل اول منتهی به 1400/03"
$ws.Cells.Item(39, 8).Value = "فصل دوم منتهی به 1400/06"
$ws.Cells.Item(39, 9).Value = "فصل سوم منتهی به 1400/09"
$ws.Cells.Item(39, 10).Value = "فصل چهارم منتهی به 1400/12"
$ws.Cells.Item(39, 11).Value = "فصل اول منتهی به 1401/03"
$ws.Cells.Item(39, 12).Value = "فصل دوم منتهی به 1401/06"
$ws.Cells.Item(39, 13).Value = "فصل سوم منتهی به 1401/09"
$ws.Cells.Item(39, 14).Value = "فصل چهارم منتهی به 1401/12"

# Row 42
$ws.Cells.Item(42, 7).Value = "-"

# Row 43
$ws.Cells.Item(43, 5).Value = 655504
$ws.Cells.Item(43, 6).Value = 358085
$ws.Cells.Item(43, 7).Value = 413800
$ws.Cells.Item(43, 8).Value = 371353
$ws.Cells.Item(43, 9).Value = 379153
$ws.Cells.Item(43, 10).Value = "-"
$ws.Cells.Item(43, 11).Value = 284520
$ws.Cells.Item(43, 12).Value = 403601
$ws.Cells.Item(43, 13).Value = 369785
$ws.Cells.Item(43, 14).Value = 706665

# Row 44
$ws.Cells.Item(44, 5).Value = 160715
$ws.Cells.Item(44, 6).Value = 32602
$ws.Cells.Item(44, 7).Value = 62079
$ws.Cells.Item(44, 8).Value = 145932
$ws.Cells.Item(44, 9).Value = 130121
$ws.Cells.Item(44, 10).Value = "-"
$ws.Cells.Item(44, 11).Value = 32179
$ws.Cells.Item(44, 12).Value = 27689
$ws.Cells.Item(44, 13).Value = 127471
$ws.Cells.Item(44, 14).Value = 53173

# Row 45
$ws.Cells.Item(45, 7).Value = "-"
$ws.Cells.Item(45, 8).Value = 0
$ws.Cells.Item(45, 10).Value = "-"
$ws.Cells.Item(45, 11).Value = 0

# Row 66
$ws.Cells.Item(66, 11).Value = 44842702
$ws.Cells.Item(66, 12).Value = 10231443
$ws.Cells.Item(66, 13).Value = 9639088
$ws.Cells.Item(66, 14).Value = 1292885

# Row 67
$ws.Cells.Item(67, 5).Value = -23826173
$ws.Cells.Item(67, 6).Value = -66017066
$ws.Cells.Item(67, 7).Value = 14980748
$ws.Cells.Item(67, 8).Value = "-"

# Row 68
$ws.Cells.Item(68, 8).Value = 46270073
$ws.Cells.Item(68, 9).Value = 23055221
$ws.Cells.Item(68, 10).Value = "-"

# Row 69
$ws.Cells.Item(69, 5).Value = -23009954
$ws.Cells.Item(69, 6).Value = -65626379
$ws.Cells.Item(69, 7).Value = 15456627
$ws.Cells.Item(69, 8).Value = 46787358
$ws.Cells.Item(69, 9).Value = 23564495
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 11).Value = 45159401
$ws.Cells.Item(69, 12).Value = 10662733
$ws.Cells.Item(69, 13).Value = 10136344
$ws.Cells.Item(69, 14).Value = 2052723

# Row 73
$ws.Cells.Item(73, 5).Value = "فصل سوم منتهی به 1399/09"
$ws.Cells.Item(73, 6).Value = "فصل چهارم منتهی به 1399/12"
$ws.Cells.Item(73, 7).Value = "فصل اول منتهی به 1400/03"
$ws.Cells.Item(73, 8).Value = "فصل دوم منتهی به 1400/06"
$ws.Cells.Item(73, 9).Value = "فصل سوم منتهی به 1400/09"
$ws.Cells.Item(73, 10).Value = "فصل چهارم منتهی به 1400/12"
$ws.Cells.Item(73, 11).Value = "فصل اول منتهی به 1401/03"
$ws.Cells.Item(73, 12).Value = "فصل دوم منتهی به 1401/06"
$ws.Cells.Item(73, 13).Value = "فصل سوم منتهی به 1401/09"
$ws.Cells.Item(73, 14).Value = "فصل چهارم منتهی به 1401/12"

# Row 76
$ws.Cells.Item(76, 7).Value = "-"

# Row 77
$ws.Cells.Item(77, 5).Value = 1408312
$ws.Cells.Item(77, 6).Value = 988963
$ws.Cells.Item(77, 7).Value = 1477384
$ws.Cells.Item(77, 8).Value = 1376650
$ws.Cells.Item(77, 9).Value = 1376148
$ws.Cells.Item(77, 10).Value = "-"
$ws.Cells.Item(77, 11).Value = 1252898
$ws.Cells.Item(77, 12).Value = 2621527
$ws.Cells.Item(77, 13).Value = 1811129
$ws.Cells.Item(77, 14).Value = 3587536

# Row 78
$ws.Cells.Item(78, 5).Value = 860105
$ws.Cells.Item(78, 6).Value = 307056
$ws.Cells.Item(78, 7).Value = 717922
$ws.Cells.Item(78, 8).Value = 1131260
$ws.Cells.Item(78, 9).Value = 407987
$ws.Cells.Item(78, 10).Value = "-"
$ws.Cells.Item(78, 11).Value = 896874
$ws.Cells.Item(78, 12).Value = 329520
$ws.Cells.Item(78, 13).Value = 412371
$ws.Cells.Item(78, 14).Value = 658009

# Row 79
$ws.Cells.Item(79, 7).Value = "-"
$ws.Cells.Item(79, 8).Value = 0
$ws.Cells.Item(79, 10).Value = "-"
$ws.Cells.Item(79, 11).Value = 0

# Row 100
$ws.Cells.Item(100, 5).Value = 341101
$ws.Cells.Item(100, 6).Value = 145367
$ws.Cells.Item(100, 7).Value = 165859
$ws.Cells.Item(100, 8).Value = 351095
$ws.Cells.Item(100, 9).Value = 196286
$ws.Cells.Item(100, 10).Value = "-"
$ws.Cells.Item(100, 11).Value = 351107
$ws.Cells.Item(100, 12).Value = 495428
$ws.Cells.Item(100, 13).Value = 348065
$ws.Cells.Item(100, 14).Value = 320076

# Row 101
$ws.Cells.Item(101, 5).Value = 2609518
$ws.Cells.Item(101, 6).Value = 1441386
$ws.Cells.Item(101, 7).Value = 2361165
$ws.Cells.Item(101, 8).Value = 2859005
$ws.Cells.Item(101, 9).Value = 1980421
$ws.Cells.Item(101, 10).Value = 0
$ws.Cells.Item(101, 11).Value = 2500879
$ws.Cells.Item(101, 12).Value = 3446475
$ws.Cells.Item(101, 13).Value = 2571565
$ws.Cells.Item(101, 14).Value = 4565621

# Row 105
$ws.Cells.Item(105, 5).Value = "فصل سوم منتهی به 1399/09"
$ws.Cells.Item(105, 6).Value = "فصل چهارم منتهی به 1399/12"
$ws.Cells.Item(105, 7).Value = "فصل اول منتهی به 1400/03"
$ws.Cells.Item(105, 8).Value = "فصل دوم منتهی به 1400/06"
$ws.Cells.Item(105, 9).Value = "فصل سوم منتهی به 1400/09"
$ws.Cells.Item(105, 10).Value = "فصل چهارم منتهی به 1400/12"
$ws.Cells.Item(105, 11).Value = "فصل اول منتهی به 1401/03"
$ws.Cells.Item(105, 12).Value = "فصل دوم منتهی به 1401/06"
$ws.Cells.Item(105, 13).Value = "فصل سوم منتهی به 1401/09"
$ws.Cells.Item(105, 14).Value = "فصل چهارم منتهی به 1401/12"

# Row 109
$ws.Cells.Item(109, 5).Value = 2148442
$ws.Cells.Item(109, 6).Value = 2761811
$ws.Cells.Item(109, 7).Value = 3570285
$ws.Cells.Item(109, 8).Value = 3707120
$ws.Cells.Item(109, 9).Value = 3629532
$ws.Cells.Item(109, 10).Value = 4140555
$ws.Cells.Item(109, 11).Value = 4403550
$ws.Cells.Item(109, 12).Value = 6495343
$ws.Cells.Item(109, 13).Value = 4897789
$ws.Cells.Item(109, 14).Value = 5076714

# Row 110
$ws.Cells.Item(110, 5).Value = 5351747
$ws.Cells.Item(110, 6).Value = 9418318
$ws.Cells.Item(110, 7).Value = 11564651
$ws.Cells.Item(110, 8).Value = 7751967
$ws.Cells.Item(110, 9).Value = 3135443
$ws.Cells.Item(110, 10).Value = 3568749
$ws.Cells.Item(110, 11).Value = 27871407
$ws.Cells.Item(110, 12).Value = 11900755
$ws.Cells.Item(110, 13).Value = 3235018
$ws.Cells.Item(110, 14).Value = 12374871

# Row 131
$ws.Cells.Item(131, 5).Value = 1164633
$ws.Cells.Item(131, 6).Value = 616469
$ws.Cells.Item(131, 7).Value = 11071
$ws.Cells.Item(131, 8).Value = 7588
$ws.Cells.Item(131, 9).Value = 8514
$ws.Cells.Item(131, 10).Value = 13843
$ws.Cells.Item(131, 11).Value = 7830
$ws.Cells.Item(131, 12).Value = 48422
$ws.Cells.Item(131, 13).Value = 36110
$ws.Cells.Item(131, 14).Value = 247567

# Row 135
$ws.Cells.Item(135, 5).Value = "فصل سوم منتهی به 1399/09"
$ws.Cells.Item(135, 6).Value = "فصل چهارم منتهی به 1399/12"
$ws.Cells.Item(135, 7).Value = "فصل اول منتهی به 1400/03"
$ws.Cells.Item(135, 8).Value = "فصل دوم منتهی به 1400/06"
$ws.Cells.Item(135, 9).Value = "فصل سوم منتهی به 1400/09"
$ws.Cells.Item(135, 10).Value = "فصل چهارم منتهی به 1400/12"
$ws.Cells.Item(135, 11).Value = "فصل اول منتهی به 1401/03"
$ws.Cells.Item(135, 12).Value = "فصل دوم منتهی به 1401/06"
$ws.Cells.Item(135, 13).Value = "فصل سوم منتهی به 1401/09"
$ws.Cells.Item(135, 14).Value = "فصل چهارم منتهی به 1401/12"

# Row 137
$ws.Cells.Item(137, 7).Value = "-"

# Row 138
$ws.Cells.Item(138, 5).Value = -611514
$ws.Cells.Item(138, 6).Value = -551772
$ws.Cells.Item(138, 7).Value = -710323
$ws.Cells.Item(138, 8).Value = -522958
$ws.Cells.Item(138, 9).Value = -799714
$ws.Cells.Item(138, 10).Value = -954209
$ws.Cells.Item(138, 11).Value = -743038
$ws.Cells.Item(138, 12).Value = -1325240
$ws.Cells.Item(138, 13).Value = -1014351
$ws.Cells.Item(138, 14).Value = -1993509

# Row 139
$ws.Cells.Item(139, 5).Value = -227861
$ws.Cells.Item(139, 6).Value = -127638
$ws.Cells.Item(139, 7).Value = -204997
$ws.Cells.Item(139, 8).Value = -376848
$ws.Cells.Item(139, 9).Value = -117285
$ws.Cells.Item(139, 10).Value = -191779
$ws.Cells.Item(139, 11).Value = -232441
$ws.Cells.Item(139, 12).Value = -198857
$ws.Cells.Item(139, 13).Value = -206391
$ws.Cells.Item(139, 14).Value = -372508

# Row 140
$ws.Cells.Item(140, 7).Value = "-"
$ws.Cells.Item(140, 8).Value = 0

# Row 141
$ws.Cells.Item(141, 5).Value = -181735
$ws.Cells.Item(141, 6).Value = -84098
$ws.Cells.Item(141, 7).Value = -112643
$ws.Cells.Item(141, 8).Value = -184944
$ws.Cells.Item(141, 9).Value = -122610
$ws.Cells.Item(141, 10).Value = -111058
$ws.Cells.Item(141, 11).Value = -197677
$ws.Cells.Item(141, 12).Value = -288717
$ws.Cells.Item(141, 13).Value = -130208
$ws.Cells.Item(141, 14).Value = -265901

# Row 142
$ws.Cells.Item(142, 5).Value = -1021110
$ws.Cells.Item(142, 6).Value = -763508
$ws.Cells.Item(142, 7).Value = -1027963
$ws.Cells.Item(142, 8).Value = -1084750
$ws.Cells.Item(142, 9).Value = -1039609
$ws.Cells.Item(142, 10).Value = -1257046
$ws.Cells.Item(142, 11).Value = -1173156
$ws.Cells.Item(142, 12).Value = -1812814
$ws.Cells.Item(142, 13).Value = -1350950
$ws.Cells.Item(142, 14).Value = -2631918

# Row 146
$ws.Cells.Item(146, 5).Value = "فصل سوم منتهی به 1399/09"
$ws.Cells.Item(146, 6).Value = "فصل چهارم منتهی به 1399/12"
$ws.Cells.Item(146, 7).Value = "فصل اول منتهی به 1400/03"
$ws.Cells.Item(146, 8).Value = "فصل دوم منتهی به 1400/06"
$ws.Cells.Item(146, 9).Value = "فصل سوم منتهی به 1400/09"
$ws.Cells.Item(146, 10).Value = "فصل چهارم منتهی به 1400/12"
$ws.Cells.Item(146, 11).Value = "فصل اول منتهی به 1401/03"
$ws.Cells.Item(146, 12).Value = "فصل دوم منتهی به 1401/06"
$ws.Cells.Item(146, 13).Value = "فصل سوم منتهی به 1401/09"
$ws.Cells.Item(146, 14).Value = "فصل چهارم منتهی به 1401/12"

# Row 148
$ws.Cells.Item(148, 7).Value = "-"

# Row 149
$ws.Cells.Item(149, 5).Value = 796798
$ws.Cells.Item(149, 6).Value = 437191
$ws.Cells.Item(149, 7).Value = 767061
$ws.Cells.Item(149, 8).Value = 853692
$ws.Cells.Item(149, 9).Value = 576434
$ws.Cells.Item(149, 10).Value = 735510
$ws.Cells.Item(149, 11).Value = 509860
$ws.Cells.Item(149, 12).Value = 1296287
$ws.Cells.Item(149, 13).Value = 796778
$ws.Cells.Item(149, 14).Value = 1594027

# Row 150
$ws.Cells.Item(150, 5).Value = 632245
$ws.Cells.Item(150, 6).Value = 179418
$ws.Cells.Item(150, 7).Value = 512925
$ws.Cells.Item(150, 8).Value = 754412
$ws.Cells.Item(150, 9).Value = 290702
$ws.Cells.Item(150, 10).Value = 285973
$ws.Cells.Item(150, 11).Value = 664433
$ws.Cells.Item(150, 12).Value = 130663
$ws.Cells.Item(150, 13).Value = 205980
$ws.Cells.Item(150, 14).Value = 285501

# Row 151
$ws.Cells.Item(151, 7).Value = "-"
$ws.Cells.Item(151, 8).Value = 0

# Row 152
$ws.Cells.Item(152, 5).Value = 159365
$ws.Cells.Item(152, 6).Value = 61269
$ws.Cells.Item(152, 7).Value = 53216
$ws.Cells.Item(152, 8).Value = 166151
$ws.Cells.Item(152, 9).Value = 73676
$ws.Cells.Item(152, 10).Value = 87806
$ws.Cells.Item(152, 11).Value = 153430
$ws.Cells.Item(152, 12).Value = 206711
$ws.Cells.Item(152, 13).Value = 217857
$ws.Cells.Item(152, 14).Value = 54175

# Row 153
$ws.Cells.Item(153, 5).Value = 1588408
$ws.Cells.Item(153, 6).Value = 677878
$ws.Cells.Item(153, 7).Value = 1333202
$ws.Cells.Item(153, 8).Value = 1774255
$ws.Cells.Item(153, 9).Value = 940812
$ws.Cells.Item(153, 10).Value = 1109289
$ws.Cells.Item(153, 11).Value = 1327723
$ws.Cells.Item(153, 12).Value = 1633661
$ws.Cells.Item(153, 13).Value = 1220615
$ws.Cells.Item(153, 14).Value = 1933703
